$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.76
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.38
$ws.Range("L2").Value = 4
$ws.Range("Y2").Value = 8.5
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 15
$ws.Range("AL2").Value = 29
$ws.Range("AO2").Value = 9.5
$ws.Range("AS2").Value = 81
$ws.Range("AZ2").Value = 23

# Row 4
$ws.Range("I4").Value = 5.5
$ws.Range("J4").Value = 2.02
$ws.Range("K4").Value = 2.22
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 10.8
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 3.72
$ws.Range("Q4").Value = 1.72
$ws.Range("R4").Value = 1.88
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.04
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.9
$ws.Range("W4").Value = 5.9
$ws.Range("X4").Value = 6.2
$ws.Range("Y4").Value = 6.8
$ws.Range("AA4").Value = 10
$ws.Range("AB4").Value = 19.5
$ws.Range("AC4").Value = 10.75
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 13.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 28
$ws.Range("AJ4").Value = 14.5
$ws.Range("AL4").Value = 40
$ws.Range("AP4").Value = 16
$ws.Range("AR4").Value = 50
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.87
$ws.Range("AU4").Value = 7.5
$ws.Range("AV4").Value = 65
$ws.Range("AX4").Value = 7.2
$ws.Range("AZ4").Value = 32
$ws.Range("BC4").Value = 400

# Row 5
$ws.Range("G5").Value = 2.2
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 2.8
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 2.38
$ws.Range("L5").Value = 3.25
$ws.Range("N5").Value = 17
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 5
$ws.Range("Q5").Value = 1.57
$ws.Range("R5").Value = 2.35
$ws.Range("U5").Value = 1.5
$ws.Range("V5").Value = 2.5
$ws.Range("Z5").Value = 21
$ws.Range("AC5").Value = 17
$ws.Range("AD5").Value = 7.5
$ws.Range("AH5").Value = 13
$ws.Range("AI5").Value = 17
$ws.Range("AL5").Value = 21
$ws.Range("AN5").Value = 4.5
$ws.Range("AO5").Value = 12
$ws.Range("AW5").Value = 351

# Row 8
$ws.Range("G8").Value = 1.5
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 5.9
$ws.Range("J8").Value = 1.98
$ws.Range("K8").Value = 2.35
$ws.Range("L8").Value = 5.4
$ws.Range("N8").Value = 13.3
$ws.Range("Q8").Value = 1.55
$ws.Range("R8").Value = 2.15
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.38
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 7.9
$ws.Range("Z8").Value = 11
$ws.Range("AB8").Value = 21
$ws.Range("AD8").Value = 8.25
$ws.Range("AE8").Value = 14.5
$ws.Range("AF8").Value = 55
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 45
$ws.Range("AJ8").Value = 18
$ws.Range("AK8").Value = 120
$ws.Range("AL8").Value = 55
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 6.9
$ws.Range("AP8").Value = 14
$ws.Range("AQ8").Value = 19.5
$ws.Range("AT8").Value = 3.25
$ws.Range("AX8").Value = 7.6
$ws.Range("AZ8").Value = 28
$ws.Range("BB8").Value = 150

# Row 9
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.9
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 1.91
$ws.Range("L9").Value = 4.5
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("O9").Value = 1.5
$ws.Range("P9").Value = 2.5
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 9
$ws.Range("Z9").Value = 19
$ws.Range("AA9").Value = 21
$ws.Range("AF9").Value = 67
$ws.Range("AH9").Value = 8.5
$ws.Range("AI9").Value = 17
$ws.Range("AK9").Value = 41
$ws.Range("AN9").Value = 4
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 29
$ws.Range("AR9").Value = 81
$ws.Range("AU9").Value = 9
$ws.Range("AX9").Value = 5.5
$ws.Range("AY9").Value = 23
$ws.Range("AZ9").Value = 34
$ws.Range("BA9").Value = 81
$ws.Range("BB9").Value = 126
$ws.Range("BC9").Value = 351

# Row 11
$ws.Range("G11").Value = 4.5
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 1.72
$ws.Range("J11").Value = 4.9
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 2.32
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 6.6
$ws.Range("O11").Value = 1.38
$ws.Range("P11").Value = 2.82
$ws.Range("Q11").Value = 2.1
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.6
$ws.Range("V11").Value = 1.72
$ws.Range("W11").Value = 11
$ws.Range("X11").Value = 25
$ws.Range("Y11").Value = 15.5
$ws.Range("Z11").Value = 80
$ws.Range("AA11").Value = 50
$ws.Range("AB11").Value = 60
$ws.Range("AC11").Value = 6.6
$ws.Range("AD11").Value = 6.7
$ws.Range("AE11").Value = 18
$ws.Range("AF11").Value = 100
$ws.Range("AG11").Value = 900
$ws.Range("AI11").Value = 7.3
$ws.Range("AJ11").Value = 8.5
$ws.Range("AK11").Value = 13
$ws.Range("AL11").Value = 15.5
$ws.Range("AN11").Value = 6.2
$ws.Range("AO11").Value = 27
$ws.Range("AP11").Value = 35
$ws.Range("AQ11").Value = 175
$ws.Range("AT11").Value = 2.6
$ws.Range("AU11").Value = 7.9
$ws.Range("AV11").Value = 80
$ws.Range("AX11").Value = 3.45
$ws.Range("AY11").Value = 8.5
$ws.Range("AZ11").Value = 20
$ws.Range("BA11").Value = 30
$ws.Range("BB11").Value = 70
